# Update code tinh luong % format cac bang
#
# Sheet "Đơn sale chính" and "Đơn phụ phẫu 1": a bunch of now-unused
# columns are removed (the report was trimmed down to the columns that
# are actually used), shifting the remaining columns left.
# Sheet "Lương": the underlying % rates changed, so several derived
# salary figures are updated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Đơn sale chính"
# Keep: Tiền tố, Mã dịch vụ, Ngày thực hiện, Cơ sở, Khách hàng,
#       Nguồn khách, Tên dịch vụ, Đơn giá gốc, Sale phụ, Upsale,
#       Đơn giá, Đã thanh toán, Tỉ lệ chiết khấu sale chính,
#       Chiết khấu sale chính
# Drop: Nhóm dịch vụ, Sale chính, Thanh toán lần đầu, Trả sau, Dư nợ,
#       Bác sĩ 1/2, Phụ phẫu 1/2, Công phụ phẫu 1/2,
#       Tỉ lệ chiết khấu sale phụ, Chiết khấu sale phụ
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Đơn sale chính")
$dropCols1 = @("AA","Y","W","V","U","T","S","R","Q","O","N","I","G")
foreach ($col in $dropCols1) {
    $ws1.Columns($col).Delete()
}

# ---------------------------------------------------------------------
# Sheet 2: "Đơn phụ phẫu 1"
# Keep: Tiền tố, Mã dịch vụ, Ngày thực hiện, Cơ sở, Khách hàng,
#       Nguồn khách, Tên dịch vụ, Phụ phẫu 1, Công phụ phẫu 1
# Drop everything else.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Đơn phụ phẫu 1")
$dropCols2 = @("AA","Z","Y","X","W","U","S","R","Q","P","O","N","M","L","K","J","I","G")
foreach ($col in $dropCols2) {
    $ws2.Columns($col).Delete()
}

# ---------------------------------------------------------------------
# Sheet 3: "Lương" - updated % rates -> recomputed figures
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Lương")
$ws3.Range("A1").Value = "Danh mục lương"
$ws3.Range("B2").Value = 17
$ws3.Range("B3").Value = 595000
$ws3.Range("B4").Value = 1821428.571428571
$ws3.Range("B12").Value = 1214285.714285714
$ws3.Range("B20").Value = 1821428.571428571
$ws3.Range("B28").Value = -1123571.428571429
$ws3.Range("B29").Value = 1214285.714285714
$ws3.Range("B30").Value = 1821428.571428571
$ws3.Range("B31").Value = 1912142.857142857
